# Working on RPM calculation.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("List1")
$ws2 = $wb.Worksheets.Item("List2")

# --- List1 (sheet1) data updates ---------------------------------------
# K2: pole pairs 6 -> 7 (drives H2 = (2*PI()*K2*G2)/60 and downstream cells)
$ws1.Range("K2").Value = 7
# F5: 9 -> 7 (drives G5/H5/I5)
$ws1.Range("F5").Value = 7

# --- List2 (sheet2): new Kth / AD / RPM table ---------------------------
# Headers - written in this order so the new shared strings land in the
# same slots as the source workbook: Kth(29), RPM(30), AD(31)
$ws2.Range("A1").Value = "Kth"
$ws2.Range("C1").Value = "RPM"
$ws2.Range("B1").Value = "AD"

$ws2.Range("A2").Value = 0.6
$ws2.Range("A3").Value = 0.7
$ws2.Range("A4").Value = 0.8
$ws2.Range("A5").Value = 0.9
$ws2.Range("A6").Value = 1
$ws2.Range("A7").Value = 1.1
$ws2.Range("A8").Value = 1.2

$ws2.Range("B2").Value = 45
$ws2.Range("B3").Value = 38
$ws2.Range("B4").Value = 33
$ws2.Range("B5").Value = 29
$ws2.Range("B6").Value = 26
$ws2.Range("B7").Value = 24
$ws2.Range("B8").Value = 22

$ws2.Range("C2").Formula = "=(1/(B2*0.0000625))*60"
$ws2.Range("C3:C8").Formula = "=(1/(B3*0.0000625))*60"

$ws2.Range("D3").Formula = "=C3-C2"
$ws2.Range("D4:D8").Formula = "=C4-C3"

# --- Chart 1 on List2: RPM vs Kth, linear trendline ----------------------
$chartObj1 = $ws2.ChartObjects().Add(788.1875, 21.75, 753.6875, 312.75)
$chartObj1.Name = "Grafikon 1"
$chart1 = $chartObj1.Chart
$chart1.ChartType = 74  # xlXYScatter
[void]$chart1.SeriesCollection().NewSeries()
$ser1 = $chart1.SeriesCollection(1)
$ser1.Name = "=List2!`$C`$1"
$ser1.XValues = $ws2.Range("A2:A8")
$ser1.Values = $ws2.Range("C2:C8")
$ser1.MarkerStyle = -4142  # xlMarkerStyleNone
$tl1 = $ser1.Trendlines().Add()
$tl1.Type = 1  # xlLinear
$tl1.DisplayEquation = $true
$tl1.DisplayRSquared = $false
$chart1.HasLegend = $true
$chart1.Legend.Position = -4152  # xlLegendPositionRight

# --- Chart 2 on List2: AD vs Kth, 3rd order polynomial trendline --------
$chartObj2 = $ws2.ChartObjects().Add(176.0625, 199.5, 753.5625, 293.25)
$chartObj2.Name = "Grafikon 2"
$chart2 = $chartObj2.Chart
$chart2.ChartType = 74  # xlXYScatter
[void]$chart2.SeriesCollection().NewSeries()
$ser2 = $chart2.SeriesCollection(1)
$ser2.Name = "=List2!`$B`$1"
$ser2.XValues = $ws2.Range("A2:A8")
$ser2.Values = $ws2.Range("B2:B8")
$ser2.MarkerStyle = -4142  # xlMarkerStyleNone
$tl2 = $ser2.Trendlines().Add()
$tl2.Type = 3  # xlPolynomial
$tl2.Order = 3
$tl2.DisplayEquation = $true
$tl2.DisplayRSquared = $false
$chart2.HasLegend = $true
$chart2.Legend.Position = -4152  # xlLegendPositionRight

# --- Selections / active sheet -------------------------------------------
[void]$ws1.Range("F6").Select()
$ws2.Activate()
[void]$ws2.Range("B8").Select()
